# Update the Publons "Test Cases" sheet (first sheet in the workbook) with
# the new/changed test cases described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 30: fix typo + change status text (blocked -> locked) -------------
$ws.Range("C30").Value = "Verify user cannot log in and it should display appropriate error message when Matching account exist in platform but in locked state."

# --- Row 40: text shifts to the "registered but not activated" message -----
$ws.Range("C40").Value = "Verify user cannot log in and it should display appropriate error message when Matching account  registered but not activated."

# --- Row 41: becomes the first of the new "facebook wrong password" rows ---
$ws.Range("B41").Value = "OPQA-5908"

# Give C41 the plain (no fill / no special font) left-aligned wrap-text style
# used by the new rows, instead of the highlighted style it had before.
$ws.Range("C34").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").HorizontalAlignment = -4131
$ws.Range("C41").Value = "Verify error message user didn't enter correct password when matching account exist in STeAM and user sign in using facebook  from login page and enters incorrect password for matching account."

$ws.Rows.Item(41).RowHeight = 30

# --- Rows 42 & 43: brand-new test cases -------------------------------------
# Start from row 40's formatting (plain bordered "text" style) for A/B/D/E...
$ws.Range("A40:E40").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)

$ws.Range("A40:E40").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)

# ...then fix up column C to use the plain left-aligned wrap-text style.
$ws.Range("C34").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C42").HorizontalAlignment = -4131

$ws.Range("C34").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("C43").HorizontalAlignment = -4131

$ws.Range("A42").Value = "PUBLONS042"
$ws.Range("B42").Value = "OPQA-5909"
$ws.Range("C42").Value = "Verify error message user didn't enter correct password when matching account exist in STeAM and user sign in using facebook  from login page and enters incorrect password for matching account."
$ws.Range("D42").Value = "Y"

$ws.Range("A43").Value = "PUBLONS043"
$ws.Range("B43").Value = "OPQA-5910"
$ws.Range("C43").Value = "Verify error message user didn't enter correct password when matching account exist in STeAM and user sign in using facebook  from login page and enters incorrect password for matching account."
$ws.Range("D43").Value = "Y"

$ws.Rows.Item(42).RowHeight = 30
$ws.Rows.Item(43).RowHeight = 30

# --- Hyperlink for the new OPQA-5908/5909/5910 block (mirrors B40:B41) -----
$ws.Hyperlinks.Add($ws.Range("B42:B43"), "https://jira.clarivate.io/browse/OPQA-5898", "", "", "https://jira.clarivate.io/browse/OPQA-5898")

# --- Cursor ends up on C37 --------------------------------------------------
$ws.Range("C37").Select()
